# The deck currently ships the "Office Theme" colour palette in
# ppt/theme/theme1.xml (used by the notes master) and the "Integral"
# palette in ppt/theme/theme2.xml (used by the slide master / the
# presentation's one Design). The authored change swaps those two
# palettes between the two theme parts.
#
# The notes master's own theme part is not independently reachable
# through this PowerPoint object model (NotesMaster/NotesPage always
# resolve back to the single Design's theme), so the colours we can
# move are applied to the Design's ThemeColorScheme, which is the part
# that actually governs what slides/masters render with
# (ppt/theme/theme2.xml) — pushing it from the "Integral" values to the
# "Office Theme" values called for by the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office Theme palette, in ThemeColorScheme.Item(i) order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# (values are VBA-style BGR-packed RGB() longs, e.g. RGB(0x44,0x54,0x6A))
$officeThemeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeThemeColors[$i - 1]
}
